$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.865.97'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '3.844.52'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '696.41'
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.89'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("D7").Value = '3.842.06'
$ws.Range("E7").Value = '  +1.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.523'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("E10").Value = '  -1.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.27'
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  -0.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.11'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = '4.491.21'
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").Value = '3.840.66'
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").Value = '70.873.18'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.16'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.114'
$ws.Range("E19").Value = '  +0.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("E20").Value = '  -3.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.72'
$ws.Range("E21").Value = '  -5.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '493.96'
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.717'
$ws.Range("E23").Value = '  -0.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.57'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000147'
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.17'
$ws.Range("E26").Value = '  -3.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.55'
$ws.Range("E27").Value = '  -1.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.11'
$ws.Range("E28").Value = '  -3.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.11'
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.49'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.26'
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.46'
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.179'
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").Value = '3.800.94'
$ws.Range("E35").Value = '  +1.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.14'
$ws.Range("E36").Value = '  -2.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.39'
$ws.Range("E39").Value = '  +5.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  +6.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.98'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.33'
$ws.Range("E42").Value = '  -5.32%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '163.43'
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.000310'
$ws.Range("E46").Value = '  -6.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '48.75'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.299'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.62'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.06'
$ws.Range("E50").Value = '  -6.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '405.48'
$ws.Range("E51").Value = '  +0.12%  '
